# participants_data.xlsx update
# - Row 4 participant replaced: "Priyanshi mandloi" / priyanshimandloi06@gmail.com -> "Pawan kushwaha" / pawankushwaha91719171@gmail.com
#   (Course stays "Designer")
# - Row 5 participant ("Asmi" / asmich1906@gmail.com / "Designer" / date) removed -> cells cleared, row kept
# - Hyperlink list: B4 now points at the new email, B5's hyperlink is gone; B2/B3 untouched

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update row 4 with the new participant ---
$ws.Range("A4").Value = "Pawan kushwaha"
$ws.Range("B4").Value = "pawankushwaha91719171@gmail.com"
$ws.Range("C4").Value = "Designer"

# --- 2. Clear row 5 contents (participant removed, row/styles remain) ---
$ws.Range("A5:D5").ClearContents()

# --- 3. Rebuild hyperlinks ---
# Stash the current formatting of a normal-hyperlink cell (B2, style used by B2/B3)
# and of a "last row" hyperlink cell (B5, style used by B4/B5) in spare cells so we
# can restore the original cell styles after re-adding the links (Hyperlinks.Add
# always reformats its target range with a fresh Hyperlink style).
$ws.Range("B2").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$ws.Range("B5").Copy()
$ws.Range("H2").PasteSpecial(-4122)

# This engine only supports clearing the *entire* hyperlink collection at once,
# so remove them all and re-add exactly the ones that should remain.
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:prince960876@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:prince960876@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B4"), "mailto:pawankushwaha91719171@gmail.com")

# Restore original cell styles that Hyperlinks.Add overwrote
$ws.Range("H1").Copy()
$ws.Range("B2").PasteSpecial(-4122)
$ws.Range("B3").PasteSpecial(-4122)
$ws.Range("H2").Copy()
$ws.Range("B4").PasteSpecial(-4122)

# Clean up the stash cells
$ws.Range("H1:H2").Clear()

# --- 4. Leave the selection where the author left it ---
[void]$ws.Range("C17").Select()
